$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.323.69'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').Value = '2.906.27'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '348.30'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.93'
$ws.Range('E6').Value = '  -6.32%  '
$ws.Range('E7').Value = '  -1.30%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  -2.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.32'
$ws.Range('E10').Value = '  -5.07%  '
$ws.Range('E11').Value = '  +1.15%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0841'
$ws.Range('E12').Value = '  -4.04%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.81'
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').Value = '3.362.43'
$ws.Range('E14').Value = '  +0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.52'
$ws.Range('E15').Value = '  -2.37%  '
$ws.Range('D16').Value = '2.908.25'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.957'
$ws.Range('E17').Value = '  -2.35%  '
$ws.Range('D18').Value = '51.305.67'
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.41'
$ws.Range('E19').Value = '  +3.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.30'
$ws.Range('E20').Value = '  -3.51%  '
$ws.Range('E21').Value = '  -5.53%  '
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.56'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '260.10'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.67'
$ws.Range('E25').Value = '  -4.13%  '
$ws.Range('E26').Value = '  +9.55%  '
$ws.Range('E27').Value = '  -4.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '26.24'
$ws.Range('E28').Value = '  -1.72%  '
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.14'
$ws.Range('E31').Value = '  -4.07%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.04'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.15'
$ws.Range('E33').Value = '  +2.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '35.26'
$ws.Range('E34').Value = '  -4.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '50.61'
$ws.Range('E35').Value = '  -4.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -6.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.08'
$ws.Range('E38').Value = '  -7.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.47'
$ws.Range('E39').Value = '  -6.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.62'
$ws.Range('E41').Value = '  -2.30%  '
$ws.Range('E42').Value = '  -1.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '22.11'
$ws.Range('E43').Value = '  -3.94%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '118.90'
$ws.Range('E44').Value = '  +7.10%  '
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('D46').Value = '2.081.33'
$ws.Range('E46').Value = '  -4.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.28'
$ws.Range('E47').Value = '  -5.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.27'
$ws.Range('E48').Value = '  -9.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.235'
$ws.Range('E49').Value = '  -5.98%  '
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.884'
$ws.Range('E51').Value = '  -5.94%  '
